# Apply edits to testdata.xlsx per commit diff ("maga 2nd time push for programming")
$wb = $excel.ActiveWorkbook

$wsOrg = $wb.Worksheets.Item("Org")
$wsContact = $wb.Worksheets.Item("contact")

# --- Update sheet1 ("Org") data ---
# Row 2: B2 changes from "POMEndtoEnd_" to "Rcb"
$wsOrg.Range("B2").Value = "Rcb"

# New rows 3-5 of data
$wsOrg.Range("B3").Value = "BengaluruBulls"
$wsOrg.Range("C3").Value = "Banking"

$wsOrg.Range("B4").Value = "BFC"
$wsOrg.Range("C4").Value = "Education"

$wsOrg.Range("B5").Value = "Karnataka"
$wsOrg.Range("C5").Value = "Finance"

# --- Selection / active sheet changes ---
# Keep "contact" sheet's own selection anchored at B6 (unchanged by the edit)
$wsContact.Range("B6").Select()

# "Org" becomes the active/displayed sheet, with its selection moved to C5
$wsOrg.Activate()
$wsOrg.Range("C5").Select()
